$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eventi")

# Insert a new row at position 37; existing rows 37-59 shift down to 38-60.
$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value = "Mostre"
$ws.Range("B37").Value = "Modena"
$ws.Range("C37").Value = "Largo Porta Sant’Agostino, 228"
$ws.Range("D37").Value = "2022-06-04T09:45:49+00:00"
$ws.Range("E37").Value = "mostra fotografica di Francesco Jodice"
$ws.Range("F37").Value = "2022-06-04T09:46:06+00:00"
$ws.Range("G37").Value = "info@agomodena.it"
$ws.Range("H37").Value = "2022-06-11T09:00:00+00:00"
$ws.Range("I37").Value = "2022-08-28T10:00:00+00:00"
$ws.Range("J37").Value = "https://www.comune.modena.it/api/novita/eventi/2022/ritratti-di-classe/@@images/3dba106b-6f95-4190-991d-b13abf85501a.jpeg"
$ws.Range("K37").Value = "Ritratti di classe"
$ws.Range("L37").Value = "2022-06-04T09:47:46+00:00"
$ws.Range("M37").Value = "AGO Modena Fabbriche culturali"
$ws.Range("N37").Value = " vedi sul sito dell'evento"
$ws.Range("O37").Value = ""
$ws.Range("P37").Value = ""
$ws.Range("Q37").Value = ""
$ws.Range("R37").Value = ""
$ws.Range("S37").Value = "Ritratti di Classe"
$ws.Range("T37").Value = ""
$ws.Range("U37").Value = "www.agomodena.it"
$ws.Range("V37").Value = $false
$ws.Range("W37").Value = 41123
$ws.Range("X37").Value = "https://www.comune.modena.it/novita/eventi/2022/ritratti-di-classe"
$ws.Range("Y37").Value = "44,64582"
$ws.Range("Z37").Value = "10,92572"
$ws.Range("AA37").Value = "POINT (10.92572 44.64582)"
